# Update crypto price/volume table to the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Leading apostrophe forces Excel to store the value as literal text
    # (the source values are plain-text, e.g. "592.08" / "67.433.46", and
    # must not be auto-coerced into numbers). Resetting the style afterwards
    # avoids leaving a stray quote-prefix / text-format style on the cell.
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "67.433.46"
Set-TextCell $ws.Range("E2") "  +1.48%  "

Set-TextCell $ws.Range("D3") "2.527.94"
Set-TextCell $ws.Range("E3") "  -1.82%  "

Set-TextCell $ws.Range("E4") "  +0.00%  "

Set-TextCell $ws.Range("D5") "592.08"
Set-TextCell $ws.Range("E5") "  +1.55%  "

Set-TextCell $ws.Range("D6") "174.75"
Set-TextCell $ws.Range("E6") "  +4.81%  "

Set-TextCell $ws.Range("E7") "  +0.03%  "

Set-TextCell $ws.Range("D8") "0.529"
Set-TextCell $ws.Range("E8") "  +0.40%  "

Set-TextCell $ws.Range("D9") "2.527.42"
Set-TextCell $ws.Range("E9") "  -1.84%  "

Set-TextCell $ws.Range("E10") "  +1.53%  "

Set-TextCell $ws.Range("E11") "  +2.41%  "

Set-TextCell $ws.Range("D12") "5.17"
Set-TextCell $ws.Range("E12") "  +0.27%  "

Set-TextCell $ws.Range("D13") "0.343"
Set-TextCell $ws.Range("E13") "  -3.16%  "

Set-TextCell $ws.Range("D14") "26.69"
Set-TextCell $ws.Range("E14") "  +0.05%  "

Set-TextCell $ws.Range("D15") "2.990.63"
Set-TextCell $ws.Range("E15") "  -1.61%  "

Set-TextCell $ws.Range("E16") "  +0.23%  "

Set-TextCell $ws.Range("D17") "67.299.65"
Set-TextCell $ws.Range("E17") "  +1.45%  "

Set-TextCell $ws.Range("D18") "2.531.69"
Set-TextCell $ws.Range("E18") "  -1.85%  "

Set-TextCell $ws.Range("E19") "  +4.98%  "

Set-TextCell $ws.Range("D20") "11.43"
Set-TextCell $ws.Range("E20") "  +0.37%  "

Set-TextCell $ws.Range("D21") "355.93"
Set-TextCell $ws.Range("E21") "  +1.66%  "

Set-TextCell $ws.Range("E22") "  -1.06%  "

Set-TextCell $ws.Range("D23") "4.64"
Set-TextCell $ws.Range("E23") "  +0.87%  "

Set-TextCell $ws.Range("E24") "  +6.68%  "

Set-TextCell $ws.Range("E25") "  +0.00%  "

Set-TextCell $ws.Range("B26") "Aptos"
Set-TextCell $ws.Range("C26") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws.Range("D26") "10.16"
Set-TextCell $ws.Range("E26") "  +2.46%  "

Set-TextCell $ws.Range("B27") "Litecoin"
Set-TextCell $ws.Range("C27") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell $ws.Range("D27") "69.78"
Set-TextCell $ws.Range("E27") "  +1.00%  "

Set-TextCell $ws.Range("D28") "1.01"
Set-TextCell $ws.Range("E28") "  +0.64%  "

Set-TextCell $ws.Range("D29") "2.658.29"
Set-TextCell $ws.Range("E29") "  -2.01%  "

Set-TextCell $ws.Range("D30") "0.0₃0986"
Set-TextCell $ws.Range("E30") "  +0.20%  "

Set-TextCell $ws.Range("D31") "557.17"
Set-TextCell $ws.Range("E31") "  +5.64%  "

Set-TextCell $ws.Range("D32") "8.22"
Set-TextCell $ws.Range("E32") "  +1.30%  "

Set-TextCell $ws.Range("E33") "  +1.65%  "

Set-TextCell $ws.Range("D35") "0.131"
Set-TextCell $ws.Range("E35") "  -0.69%  "

Set-TextCell $ws.Range("E36") "  -0.01%  "

Set-TextCell $ws.Range("E37") "  +1.20%  "

Set-TextCell $ws.Range("D38") "157.33"
Set-TextCell $ws.Range("E38") "  +0.58%  "

Set-TextCell $ws.Range("D39") "18.71"
Set-TextCell $ws.Range("E39") "  -0.23%  "

Set-TextCell $ws.Range("D40") "18.46"
Set-TextCell $ws.Range("E40") "  +0.86%  "

Set-TextCell $ws.Range("E41") "  -0.96%  "

Set-TextCell $ws.Range("E42") "  +2.13%  "

Set-TextCell $ws.Range("D43") "5.15"
Set-TextCell $ws.Range("E43") "  +0.79%  "

Set-TextCell $ws.Range("D44") "2.55"
Set-TextCell $ws.Range("E44") "  +5.07%  "

Set-TextCell $ws.Range("E45") "  -0.01%  "

Set-TextCell $ws.Range("D46") "149.26"
Set-TextCell $ws.Range("E46") "  +0.10%  "

Set-TextCell $ws.Range("D47") "0.560"
Set-TextCell $ws.Range("E47") "  -0.90%  "

Set-TextCell $ws.Range("D48") "0.0₆0277"
Set-TextCell $ws.Range("E48") "  -3.43%  "

Set-TextCell $ws.Range("D49") "3.69"
Set-TextCell $ws.Range("E49") "  -0.57%  "

Set-TextCell $ws.Range("E50") "  -1.02%  "

Set-TextCell $ws.Range("E51") "  -0.18%  "
